$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '27.680.97'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.98%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.621.89'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.85%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.993'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.82%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '210.05'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.05%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.517'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.31%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.992'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.77%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '23.21'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.54%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.256'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.91%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0607'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.35%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0876'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.86%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.852.03'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.84%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.641.61'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.33%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '3.99'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.73%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.559'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.22%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '64.79'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.98%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '27.711.52'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.89%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '227.75'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.42%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '7.68'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.63%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.0₃0716'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.09%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.994'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.65%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.32'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.26%  '

# Row 23
$ws.Range("E23").Value = '  -2.79%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.05'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.84%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '154.64'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.15%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '6.92'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.74%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.110'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.89%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '15.46'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.35%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.994'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.70%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.17'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -1.14%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.0478'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.07%  '

# Row 32
$ws.Range("E32").Value = '  -0.35%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.08'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.11%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.391.68'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.15%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.59'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.52%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.24%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.33'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -1.37%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.0170'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.33%  '

# Row 39
$ws.Range("E39").Value = '  -1.09%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.846'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -3.13%  '

# Row 41
$ws.Range("E41").Value = '  -0.99%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.994'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.63%  '

# Row 43
$ws.Range("E43").Value = '  -0.33%  '

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '65.60'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.95%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '5.39'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.69%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.767.04'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.58%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.16'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.83%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '87.91'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.01%  '

# Row 49
$ws.Range("E49").Value = '  +1.28%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0503'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.66%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '7.56'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.96%  '
